$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 905, shifting existing rows 905:981 down to 906:982
$ws.Rows("905:905").Insert()

# Populate the new row 905 with a fresh "1a plateado" Limón price record.
# Columns A,B,C,E,F,G,H,I,J,K,L,R keep the same values the old row 905 had;
# D,M,N,O,P,Q,S,T get new values.
$ws.Range("A905").Value = 7
$ws.Range("B905").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C905").Value = "Ñuble"
$ws.Range("D905").Value = 45013
$ws.Range("E905").Value = 16
$ws.Range("F905").Value = "Fruta"
$ws.Range("G905").Value = 100102
$ws.Range("H905").Value = "Cítricos"
$ws.Range("I905").Value = 100102003
$ws.Range("J905").Value = "Limón"
$ws.Range("K905").Value = "Sin especificar"
$ws.Range("L905").Value = "1a plateado"
$ws.Range("M905").Value = 250
$ws.Range("N905").Value = 25000
$ws.Range("O905").Value = 26000
$ws.Range("P905").Value = 25400
$ws.Range("Q905").Value = "`$/malla 18 kilos"
$ws.Range("R905").Value = "Región de O'Higgins"
$ws.Range("S905").Value = 1411
$ws.Range("T905").Value = 18
